# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp in A1
# - Update several countries' case statistics
# - Two pairs of countries swap their relative order (new data pushes one
#   country above its former neighbour); the neighbour's row keeps its old,
#   unchanged figures but moves down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 10:19"

# --- Simple numeric refreshes (country stays on the same row) -----------

# Rusia (row 7)
$ws.Cells.Item(7, 2).Value = 1225889
$ws.Cells.Item(7, 3).Value = 10888
$ws.Cells.Item(7, 4).Value = 982324
$ws.Cells.Item(7, 5).Value = 222090
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 117
$ws.Cells.Item(7, 8).Value = 21475

# Hungria (row 77)
$ws.Cells.Item(77, 2).Value = 31480
$ws.Cells.Item(77, 3).Value = 905
$ws.Cells.Item(77, 4).Value = 8165
$ws.Cells.Item(77, 5).Value = 22482
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 11
$ws.Cells.Item(77, 8).Value = 833

# Eslovaquia (row 98)
$ws.Cells.Item(98, 2).Value = 13492
$ws.Cells.Item(98, 3).Value = 353
$ws.Cells.Item(98, 4).Value = 4865
$ws.Cells.Item(98, 5).Value = 8572
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 55

# Estonia (row 140)
$ws.Cells.Item(140, 2).Value = 3617
$ws.Cells.Item(140, 3).Value = 10
$ws.Cells.Item(140, 4).Value = 2755
$ws.Cells.Item(140, 5).Value = 795
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 67

# --- Pairs that swap order: updated country moves to the upper row, ------
# --- previous occupant keeps its old values and drops to the lower row ---

# Filipinas now ranks above Turquia (rows 22-23)
$ws.Cells.Item(22, 1).Value = "Filipinas"
$ws.Cells.Item(22, 2).Value = 324762
$ws.Cells.Item(22, 3).Value = 2291
$ws.Cells.Item(22, 4).Value = 273123
$ws.Cells.Item(22, 5).Value = 45799
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 64
$ws.Cells.Item(22, 8).Value = 5840

$ws.Cells.Item(23, 1).Value = "Turquia"
$ws.Cells.Item(23, 2).Value = 324443
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 285050
$ws.Cells.Item(23, 5).Value = 30952
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 8441

# Croacia now ranks above Birmania (rows 89-90)
$ws.Cells.Item(89, 1).Value = "Croacia"
$ws.Cells.Item(89, 2).Value = 17797
$ws.Cells.Item(89, 3).Value = 138
$ws.Cells.Item(89, 4).Value = 16031
$ws.Cells.Item(89, 5).Value = 1466
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 300

$ws.Cells.Item(90, 1).Value = "Birmania"
$ws.Cells.Item(90, 2).Value = 17794
$ws.Cells.Item(90, 3).Value = 0
$ws.Cells.Item(90, 4).Value = 5195
$ws.Cells.Item(90, 5).Value = 12187
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 412

# Letonia now ranks above Principado de Andorra (rows 156-157)
$ws.Cells.Item(156, 1).Value = "Letonia"
$ws.Cells.Item(156, 2).Value = 2126
$ws.Cells.Item(156, 3).Value = 40
$ws.Cells.Item(156, 4).Value = 1307
$ws.Cells.Item(156, 5).Value = 780
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 1
$ws.Cells.Item(156, 8).Value = 39

$ws.Cells.Item(157, 1).Value = "Principado de Andorra"
$ws.Cells.Item(157, 2).Value = 2110
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 1540
$ws.Cells.Item(157, 5).Value = 517
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 53

# Islas Malvinas now ranks above Montserrat (rows 215-216)
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 6).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 1

Write-Output "edit complete"
